$d = $word.ActiveDocument

# 1. Merge "Desig" + "n, implement and deploy..." into one run of text
#    (this is purely a textual no-op from the reader's perspective - the
#    rendered text stays the same, so a Find/Replace on the full text is safe)
$d.Content.Find.Execute("Desig" + "n, implement and deploy foundational platform services in AWS for the Aura digital security tool", $true, $false, $false, $false, $false, $true, 1, $false, "Design, implement and deploy foundational platform services in AWS for the Aura digital security tool", 2) | Out-Null

# 2. Reorder / extend the sentence about the service mesh:
#    "... service mesh with Cassandra 4.0 as the data store and React 17.0 as the front end."
# becomes
#    "... service mesh with React 17.0 as the front end, Cassandra 4.0 as the data store, and Keycloak 15.0 as the identity management platform."
$d.Content.Find.Execute("service mesh with Cassandra 4.0 as the data store and React 17.0 as the front end", $true, $false, $false, $false, $false, $true, 1, $false, "service mesh with React 17.0 as the front end, Cassandra 4.0 as the data store, and Keycloak 15.0 as the identity management platform", 2) | Out-Null
